$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B13").Value = 0.052444562315940857
$ws.Range("C13").Value = 0.14833562076091766
$ws.Range("B14").Value = 0.22181642055511475
$ws.Range("C14").Value = 0.62739157676696777
$ws.Range("B15").Value = 1.0168312788009644
$ws.Range("C15").Value = 1.1753454208374023
$ws.Range("B16").Value = 3.3291506767272949
$ws.Range("C16").Value = 1.349273681640625
$ws.Range("B17").Value = 6.5806617736816406
$ws.Range("C17").Value = 1.4236965179443359
$ws.Range("B18").Value = 9.5487251281738281
$ws.Range("C18").Value = 1.8933607339859009
$ws.Range("B19").Value = 11.136787414550781
$ws.Range("C19").Value = 2.7605080604553223
$ws.Range("B20").Value = 10.8428316116333
$ws.Range("C20").Value = 3.7245430946350098
$ws.Range("B21").Value = 8.9065876007080078
$ws.Range("C21").Value = 4.2764959335327148
$ws.Range("B22").Value = 6.1420974731445313
$ws.Range("C22").Value = 3.8727495670318604
$ws.Range("B23").Value = 3.4864053726196289
$ws.Range("C23").Value = 2.4930763244628906
$ws.Range("B24").Value = 1.3679167032241821
$ws.Range("C24").Value = 1.2112330198287964
$ws.Range("B25").Value = 0.25140956044197083
$ws.Range("C25").Value = 0.36038661003112793
$ws.Range("B26").Value = 0
$ws.Range("C26").Value = 0
$ws.Range("B28").Value = 0.031898848712444305
$ws.Range("C28").Value = 0.090223565697669983
$ws.Range("B29").Value = 0.077550873160362244
$ws.Range("C29").Value = 0.21934698522090912
$ws.Range("B30").Value = 0.13196295499801636
$ws.Range("C30").Value = 0.31111562252044678
$ws.Range("B31").Value = 0.17832286655902863
$ws.Range("C31").Value = 0.34652248024940491
$ws.Range("B32").Value = 0.18989077210426331
$ws.Range("C32").Value = 0.35453689098358154
$ws.Range("B33").Value = 0.15800571441650391
$ws.Range("C33").Value = 0.33755061030387878
$ws.Range("B34").Value = 0.10072866082191467
$ws.Range("C34").Value = 0.26595574617385864
$ws.Range("B35").Value = 0.080607175827026367
$ws.Range("C35").Value = 0.15380460023880005
$ws.Range("B36").Value = 0.096128210425376892
$ws.Range("C36").Value = 0.2407984733581543
$ws.Range("B37").Value = 0.11763719469308853
$ws.Range("C37").Value = 0.33272823691368103
$ws.Range("B38").Value = 0.10575366020202637
$ws.Range("C38").Value = 0.278071790933609
$ws.Range("B39").Value = 0.077496379613876343
$ws.Range("C39").Value = 0.14555294811725616
$ws.Range("B40").Value = 0.054860301315784454
$ws.Range("C40").Value = 0.15516836941242218
$ws.Range("B41").Value = 0.051519401371479034
$ws.Range("C41").Value = 0.14571887254714966
$ws.Range("B42").Value = 0.025638947263360023
$ws.Range("C42").Value = 0.072517894208431244
$ws.Range("B43").Value = 0.0087056923657655716
$ws.Range("C43").Value = 0.021260660141706467
$ws.Range("B44").Value = 0.031609933823347092
$ws.Range("C44").Value = 0.089406393468379974
$ws.Range("B45").Value = 0.062704741954803467
$ws.Range("C45").Value = 0.1431535929441452
$ws.Range("B46").Value = 0.14257071912288666
$ws.Range("C46").Value = 0.28176695108413696
$ws.Range("B47").Value = 0.23360620439052582
$ws.Range("C47").Value = 0.58185297250747681
$ws.Range("B48").Value = 0.28969621658325195
$ws.Range("C48").Value = 0.80839776992797852
$ws.Range("B49").Value = 0.28610289096832275
$ws.Range("C49").Value = 0.80922126770019531
$ws.Range("B50").Value = 0.24020771682262421
$ws.Range("C50").Value = 0.5845034122467041
$ws.Range("B51").Value = 0.18962246179580688
$ws.Range("C51").Value = 0.35147830843925476
$ws.Range("B52").Value = 0.16739143431186676
$ws.Range("C52").Value = 0.42080292105674744
$ws.Range("B53").Value = 0.1791745126247406
$ws.Range("C53").Value = 0.50678199529647827
$ws.Range("B54").Value = 0.16317477822303772
$ws.Range("C54").Value = 0.46152797341346741
$ws.Range("B55").Value = 0.11678589880466461
$ws.Range("C55").Value = 0.30838626623153687
$ws.Range("B56").Value = 0.065147385001182556
$ws.Range("C56").Value = 0.12775759398937225
$ws.Range("B57").Value = 0.074515827000141144
$ws.Range("C57").Value = 0.16343432664871216
$ws.Range("B58").Value = 0.14588671922683716
$ws.Range("C58").Value = 0.39918354153633118
$ws.Range("B59").Value = 0.21449357271194458
$ws.Range("C59").Value = 0.60667943954467773
$ws.Range("B60").Value = 0.24202245473861694
$ws.Range("C60").Value = 0.68454283475875854
$ws.Range("B61").Value = 0.30074435472488403
$ws.Range("C61").Value = 0.61385041475296021
$ws.Range("B62").Value = 0.49727684259414673
$ws.Range("C62").Value = 1.0408883094787598
$ws.Range("B63").Value = 0.80976420640945435
$ws.Range("C63").Value = 2.12850022315979
$ws.Range("B64").Value = 1.1429497003555298
$ws.Range("C64").Value = 3.2291028499603271
$ws.Range("B65").Value = 1.3697460889816284
$ws.Range("C65").Value = 3.8742270469665527
$ws.Range("B66").Value = 1.3504477739334106
$ws.Range("C66").Value = 3.8196432590484619
$ws.Range("B67").Value = 1.0819442272186279
$ws.Range("C67").Value = 3.0602004528045654
$ws.Range("B68").Value = 0.65969276428222656
$ws.Range("C68").Value = 1.8658928871154785
$ws.Range("B69").Value = 0.28349402546882629
$ws.Range("C69").Value = 0.69444751739501953
$ws.Range("B70").Value = 0.12694679200649261
$ws.Range("C70").Value = 0.32492616772651672
$ws.Range("B71").Value = 0.21474792063236237
$ws.Range("C71").Value = 0.6073988676071167
$ws.Range("B72").Value = 0.29213255643844604
$ws.Range("C72").Value = 0.826275646686554
